$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.479.27"
$ws.Range("E2").Value = "  +1.33%  "

$ws.Range("D3").Value = "3.172.41"
$ws.Range("E3").Value = "  -0.98%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'596.39"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").Value = "'135.62"
$ws.Range("E6").Value = "  -0.49%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("D8").Value = "3.168.32"
$ws.Range("E8").Value = "  -1.07%  "

$ws.Range("E9").Value = "  +1.70%  "

$ws.Range("E10").Value = "  -1.27%  "

$ws.Range("D11").Value = "'5.33"
$ws.Range("E11").Value = "  -0.51%  "

$ws.Range("D12").Value = "'0.455"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("E13").Value = "  -0.07%  "

$ws.Range("D14").Value = "'34.75"
$ws.Range("E14").Value = "  +3.12%  "

$ws.Range("D15").Value = "3.693.90"
$ws.Range("E15").Value = "  -1.00%  "

$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").Value = "3.172.93"
$ws.Range("E17").Value = "  -1.54%  "

$ws.Range("D18").Value = "63.488.89"
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("E19").Value = "  -2.39%  "

$ws.Range("D20").Value = "'461.99"
$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("D21").Value = "'13.96"
$ws.Range("E21").Value = "  -0.58%  "

$ws.Range("E22").Value = "  -2.61%  "

$ws.Range("D23").Value = "'7.67"
$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'83.07"
$ws.Range("E24").Value = "  -0.72%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'13.21"
$ws.Range("E25").Value = "  -3.32%  "

$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("E27").Value = "  -1.95%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("D30").Value = "'7.73"
$ws.Range("E30").Value = "  -3.03%  "

$ws.Range("D31").Value = "'6.80"
$ws.Range("E31").Value = "  -2.00%  "

$ws.Range("D32").Value = "'27.31"
$ws.Range("E32").Value = "  -1.17%  "

$ws.Range("E33").Value = "  -1.61%  "

$ws.Range("E34").Value = "  -1.80%  "

$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").Value = "'51.47"
$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("E38").Value = "  +4.47%  "

$ws.Range("D39").Value = "'0.0390"
$ws.Range("E39").Value = "  -1.00%  "

$ws.Range("D40").Value = "'8.12"
$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("E41").Value = "  -2.50%  "

$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").Value = "'391.93"
$ws.Range("E43").Value = "  -6.98%  "

$ws.Range("D44").Value = "2.790.03"
$ws.Range("E44").Value = "  -7.64%  "

$ws.Range("E45").Value = "  -1.47%  "

$ws.Range("D46").Value = "'128.05"
$ws.Range("E46").Value = "  +1.84%  "

$ws.Range("D47").Value = "'35.90"
$ws.Range("E47").Value = "  -1.30%  "

$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("E49").Value = "  -2.99%  "

$ws.Range("D50").Value = "'25.22"
$ws.Range("E50").Value = "  -3.63%  "

$ws.Range("E51").Value = "  -0.95%  "
